# "Generate Report for Archive"
# Localization status report refresh: the items previously queued for
# handoff have moved into active translation, and the (now shorter)
# "Status" column is narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newColumnWidth = 12.5   # yields the narrower "Status" column width

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
